$wb = $excel.ActiveWorkbook

# This script applies updated market-price derived values (scheduled runner sync)
# to the H:N "profit" columns across multiple class sheets.

$ws = $wb.Worksheets.Item("ALC")
# Row 41
$ws.Range("H41").Value = 840.43475
$ws.Range("I41").Value = 919.58826
$ws.Range("K41").Value = 919.58826
$ws.Range("M41").Value = -479.58826
# Row 64
$ws.Range("H64").Value = 6331.5454
$ws.Range("I64").Value = 3608.3333
$ws.Range("K64").Value = 3608.3333
$ws.Range("M64").Value = -3360.3333
# Row 67
$ws.Range("H67").Value = 6331.5454
$ws.Range("I67").Value = 3608.3333
$ws.Range("K67").Value = 3608.3333
$ws.Range("M67").Value = -2750.3333
# Row 137
$ws.Range("H137").Value = 1846.1864
$ws.Range("I137").Value = 1850.95
$ws.Range("J137").Value = 1836.1578
$ws.Range("K137").Value = 5552.85
$ws.Range("L137").Value = 5508.4734
$ws.Range("M137").Value = -3002.85
$ws.Range("N137").Value = -10608.4734

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 4499
$ws.Range("I61").Value = 3503.5
$ws.Range("J61").Value = 5992.25
$ws.Range("K61").Value = 3503.5
$ws.Range("L61").Value = 5992.25
$ws.Range("M61").Value = -3291.5
$ws.Range("N61").Value = -6416.25
# Row 122
$ws.Range("H122").Value = 1993
$ws.Range("I122").Value = 1994
$ws.Range("K122").Value = 5982
$ws.Range("M122").Value = -3532
# Row 132
$ws.Range("H132").Value = 4079.9167
$ws.Range("I132").Value = 3712.15
$ws.Range("K132").Value = 11136.45
$ws.Range("M132").Value = -8606.450000000001
# Row 136
$ws.Range("H136").Value = 4499
$ws.Range("I136").Value = 3503.5
$ws.Range("J136").Value = 5992.25
$ws.Range("K136").Value = 10510.5
$ws.Range("L136").Value = 17976.75
$ws.Range("M136").Value = -7960.5
$ws.Range("N136").Value = -23076.75

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 4008.7144
$ws.Range("I134").Value = 4154.4614
$ws.Range("K134").Value = 12463.3842
$ws.Range("M134").Value = -9928.3842

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 5682.25
$ws.Range("J16").Value = 6656
$ws.Range("L16").Value = 6656
$ws.Range("N16").Value = -7230
# Row 31
$ws.Range("H31").Value = 1473.25
$ws.Range("I31").Value = 1317.2122
$ws.Range("J31").Value = 3189.6667
$ws.Range("K31").Value = 1317.2122
$ws.Range("L31").Value = 3189.6667
$ws.Range("M31").Value = -1022.2122
$ws.Range("N31").Value = -3779.6667
# Row 34
$ws.Range("H34").Value = 1473.25
$ws.Range("I34").Value = 1317.2122
$ws.Range("J34").Value = 3189.6667
$ws.Range("K34").Value = 1317.2122
$ws.Range("L34").Value = 3189.6667
$ws.Range("M34").Value = -1115.2122
$ws.Range("N34").Value = -3593.6667
# Row 55
$ws.Range("H55").Value = 34693.332
$ws.Range("J55").Value = 34693.332
$ws.Range("L55").Value = 34693.332
$ws.Range("N55").Value = -35323.332
# Row 99
$ws.Range("H99").Value = 3161
$ws.Range("I99").Value = 3008.4546
$ws.Range("K99").Value = 3008.4546
$ws.Range("M99").Value = -1510.4546
# Row 113
$ws.Range("H113").Value = 5682.25
$ws.Range("J113").Value = 6656
$ws.Range("L113").Value = 6656
$ws.Range("N113").Value = -10996
# Row 122
$ws.Range("H122").Value = 3498.5715
$ws.Range("J122").Value = 5997.5
$ws.Range("L122").Value = 17992.5
$ws.Range("N122").Value = -22892.5
# Row 126
$ws.Range("H126").Value = 3161
$ws.Range("I126").Value = 3008.4546
$ws.Range("K126").Value = 9025.363799999999
$ws.Range("M126").Value = -6555.363799999999
# Row 132
$ws.Range("H132").Value = 2773.4666
$ws.Range("I132").Value = 2800.25
$ws.Range("J132").Value = 2666.3333
$ws.Range("K132").Value = 8400.75
$ws.Range("L132").Value = 7998.999899999999
$ws.Range("M132").Value = -5870.75
$ws.Range("N132").Value = -13058.9999
# Row 134
$ws.Range("H134").Value = 12038.5
$ws.Range("I134").Value = 8598.966
$ws.Range("K134").Value = 25796.898
$ws.Range("M134").Value = -23261.898

$ws = $wb.Worksheets.Item("CUL")
# Row 80
$ws.Range("H80").Value = 4013.2666
$ws.Range("J80").Value = 4207.6924
$ws.Range("L80").Value = 12623.0772
$ws.Range("N80").Value = -14495.0772
# Row 83
$ws.Range("H83").Value = 4013.2666
$ws.Range("J83").Value = 4207.6924
$ws.Range("L83").Value = 37869.2316
$ws.Range("N83").Value = -47229.2316
# Row 92
$ws.Range("H92").Value = 1434
$ws.Range("J92").Value = 1751
$ws.Range("L92").Value = 5253
$ws.Range("N92").Value = -7749
# Row 97
$ws.Range("H97").Value = 1302.8572
$ws.Range("I97").Value = 350
$ws.Range("J97").Value = 1461.6666
$ws.Range("K97").Value = 1050
$ws.Range("L97").Value = 4384.9998
$ws.Range("M97").Value = -554
$ws.Range("N97").Value = -5376.9998
# Row 132
$ws.Range("H132").Value = 2054.1875
$ws.Range("I132").Value = 1540.2941
$ws.Range("J132").Value = 2636.6
$ws.Range("K132").Value = 13862.6469
$ws.Range("L132").Value = 23729.4
$ws.Range("M132").Value = -11332.6469
$ws.Range("N132").Value = -28789.4

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 11204.875
$ws.Range("J70").Value = 14176.454
$ws.Range("L70").Value = 14176.454
$ws.Range("N70").Value = -14716.454
# Row 73
$ws.Range("H73").Value = 11204.875
$ws.Range("J73").Value = 14176.454
$ws.Range("L73").Value = 14176.454
$ws.Range("N73").Value = -16048.454
# Row 97
$ws.Range("H97").Value = 1621.762
$ws.Range("I97").Value = 1840.1428
$ws.Range("J97").Value = 1185
$ws.Range("K97").Value = 1840.1428
$ws.Range("L97").Value = 1185
$ws.Range("M97").Value = -1344.1428
$ws.Range("N97").Value = -2177
# Row 122
$ws.Range("H122").Value = 1828.6
$ws.Range("I122").Value = 1535.75
$ws.Range("K122").Value = 4607.25
$ws.Range("M122").Value = -2157.25
# Row 132
$ws.Range("H132").Value = 5961.125
$ws.Range("I132").Value = 4737.8
$ws.Range("K132").Value = 14213.4
$ws.Range("M132").Value = -11683.4

$ws = $wb.Worksheets.Item("LTW")
# Row 29
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
# Row 82
$ws.Range("H82").Value = 17019
$ws.Range("J82").Value = 3360
$ws.Range("L82").Value = 3360
$ws.Range("N82").Value = -4082
# Row 85
$ws.Range("H85").Value = 17019
$ws.Range("J85").Value = 3360
$ws.Range("L85").Value = 3360
$ws.Range("N85").Value = -5856

$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 1479.1428
$ws.Range("I96").Value = 1212.5
$ws.Range("J96").Value = 1834.6666
$ws.Range("K96").Value = 1212.5
$ws.Range("L96").Value = 1834.6666
$ws.Range("M96").Value = 160.5
$ws.Range("N96").Value = -4580.6666
# Row 122
$ws.Range("H122").Value = 2493.4243
$ws.Range("I122").Value = 2272.32
$ws.Range("K122").Value = 6816.960000000001
$ws.Range("M122").Value = -4366.960000000001
# Row 136
$ws.Range("H136").Value = 2167.2917
$ws.Range("I136").Value = 1900.8
$ws.Range("J136").Value = 3499.75
$ws.Range("K136").Value = 5702.4
$ws.Range("L136").Value = 10499.25
$ws.Range("M136").Value = -3152.4
$ws.Range("N136").Value = -15599.25

